# update stats, fix names...
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("roster")

function Sort-RosterColumn {
    param($colLetter, $rowStart, $rowEnd)
    $rng = $ws.Range("$colLetter$rowStart" + ":" + "$colLetter$rowEnd")
    $srt = $ws.Sort
    $srt.SortFields.Clear()
    $srt.SortFields.Add($rng)
    $srt.SetRange($rng)
    $srt.Header = 0
    $srt.Apply()
}

# Each roster column is independently sorted A->Z (names), the names
# lists are of different lengths so the ranges differ per column.
Sort-RosterColumn "A" 2 13
Sort-RosterColumn "B" 2 16
Sort-RosterColumn "C" 2 12
Sort-RosterColumn "D" 2 12

# Make "roster" the active / selected sheet (was "8-2").
$ws.Activate()
